# Version 1.1 de la solapa FacturasPAMI - Agregado de parametros en la query
# sql de NuevoPami
#
# Adds a new "TipoPami" column (F) to the "Tabla2" table on Hoja1, classifying
# each "Tipo Comprobante" row as Facturado / No Aplica / Debitos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add a new table column. This extends the table (and the worksheet
# dimension/autofilter) from A1:E11 to A1:F11 automatically.
$newCol = $lo.ListColumns.Add()

# Header text for the new column
$ws.Range("F1").Value2 = "TipoPami"

# Match the look of the rest of the table: centered + wrapping header,
# centered data cells (same formatting already used by the other columns).
$ws.Range("F1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("F1").WrapText = $true
$ws.Range("F2:F11").HorizontalAlignment = -4108   # xlCenter

# Populate the new column's values
$ws.Range("F2").Value2 = "Facturado"
$ws.Range("F3").Value2 = "Facturado"
$ws.Range("F4").Value2 = "Facturado"
$ws.Range("F5").Value2 = "Facturado"
$ws.Range("F6").Value2 = "No Aplica"
$ws.Range("F7").Value2 = "No Aplica"
$ws.Range("F8").Value2 = "No Aplica"
$ws.Range("F9").Value2 = "Debitos"
$ws.Range("F10").Value2 = "Debitos"
$ws.Range("F11").Value2 = "Debitos"

# Leave the selection where the author left it
$ws.Range("F6").Select()
